$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column (D) formatted as Text so numeric-looking values
# (e.g. "0.9988", "1.000") are not reinterpreted as numbers and lose their
# original formatting (trailing zeros, exact decimal text, etc).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.189.92"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "1.644.72"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").Value = "308.47"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "0.3901"
$ws.Range("E7").Value = "  -1.08%  "
$ws.Range("D8").Value = "0.3870"
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("D9").Value = "1.002"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "49.72"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").Value = "1.354"
$ws.Range("E11").Value = "  -5.13%  "
$ws.Range("D12").Value = "0.08704"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "23.72"
$ws.Range("E13").Value = "  -5.57%  "
$ws.Range("D14").Value = "7.116"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").Value = "0.00001295"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "7.469"
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("D17").Value = "1.611.17"
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").Value = "95.37"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "0.06913"
$ws.Range("E19").Value = "  -3.12%  "
$ws.Range("D20").Value = "20.58"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").Value = "6.911"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "13.60"
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("D24").Value = "24.175.59"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").Value = "2.333"
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("D26").Value = "2.770"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "22.39"
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("D28").Value = "157.96"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("D29").Value = "8.482"
$ws.Range("E29").Value = "  +7.61%  "
$ws.Range("D30").Value = "140.50"
$ws.Range("E30").Value = "  -5.72%  "
$ws.Range("D31").Value = "5.349"
$ws.Range("E31").Value = "  -10.20%  "
$ws.Range("D32").Value = "2.419"
$ws.Range("E32").Value = "  -8.72%  "
$ws.Range("D33").Value = "1.817.71"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").Value = "6.977"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").Value = "0.08054"
$ws.Range("E35").Value = "  -4.59%  "
$ws.Range("D36").Value = "0.02906"
$ws.Range("E36").Value = "  -5.70%  "
$ws.Range("D37").Value = "0.2686"
$ws.Range("E37").Value = "  -5.00%  "
$ws.Range("D38").Value = "0.9502"
$ws.Range("E38").Value = "  -6.61%  "
$ws.Range("D39").Value = "0.09215"
$ws.Range("E39").Value = "  -3.95%  "
$ws.Range("D40").Value = "1.464"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "9.950"
$ws.Range("E41").Value = "  -6.26%  "
$ws.Range("D42").Value = "0.7574"
$ws.Range("E42").Value = "  -5.46%  "
$ws.Range("D43").Value = "13.05"
$ws.Range("E43").Value = "  -5.00%  "
$ws.Range("D44").Value = "15.93"
$ws.Range("E44").Value = "  -5.00%  "
$ws.Range("D45").Value = "0.6919"
$ws.Range("E45").Value = "  -3.89%  "
$ws.Range("D46").Value = "2.470"
$ws.Range("E46").Value = "  -4.90%  "
$ws.Range("D47").Value = "4.092"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "0.08406"
$ws.Range("E49").Value = "  -4.21%  "
$ws.Range("D50").Value = "1.262"
$ws.Range("E50").Value = "  -7.52%  "
$ws.Range("D51").Value = "133.11"
$ws.Range("E51").Value = "  -4.44%  "
